# Fruta / hortaliza, semanal
#
# A new weekly price record for Mango (Vega Modelo de Temuco) needs to be
# inserted as row 94 of the data table, pushing every subsequent record
# down by one row (old row 94 -> new row 95, ..., old row 212 -> new row 213).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 94, shifting rows 94:212 down to 95:213.
$ws.Rows("94:94").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 44467
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100108
$ws.Range("H94").Value = "Tropicales y subtropicales"
$ws.Range("I94").Value = 100108002
$ws.Range("J94").Value = "Mango"
$ws.Range("K94").Value = "Sin especificar"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 200
$ws.Range("N94").Value = 9000
$ws.Range("O94").Value = 9000
$ws.Range("P94").Value = 9000
$ws.Range("Q94").Value = "$/bandeja 4 kilos"
$ws.Range("R94").Value = "Brasil"
$ws.Range("S94").Value = 2250
$ws.Range("T94").Value = 4
